$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Actual Production (MW)" values for rows 2..39 (rows 40..97 stay at 0, unchanged)
$newB = @{
  2 = 408
  3 = 339
  4 = 289
  5 = 268
  6 = 243
  7 = 228
  8 = 205
  9 = 189
  10 = 177
  11 = 170
  12 = 173
  13 = 175
  14 = 165
  15 = 155
  16 = 149
  17 = 149
  18 = 139
  19 = 132
  20 = 132
  21 = 123
  22 = 106
  23 = 100
  24 = 98
  25 = 100
  26 = 97
  27 = 86
  28 = 82
  29 = 78
  30 = 66
  31 = 48
  32 = 79
  33 = 51
  34 = 38
  35 = 48
  36 = 68
  37 = 89
  38 = 111
  39 = 123
}

# Shift every Timestamp (column A) forward by 18 days, keeping the time-of-day fraction
for ($r = 2; $r -le 97; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $oldSerial = $cellA.Value2()
    $cellA.Value = $oldSerial + 18

    if ($newB.ContainsKey($r)) {
        $ws.Cells.Item($r, 2).Value = $newB[$r]
    }
}
